$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on each touched Price/Volume cell individually so that
# numeric-looking strings (e.g. "0.999", "143.80") are kept as literal text,
# matching the original inlineStr cell content (no numeric coercion, no loss
# of trailing zeros). NumberFormat is applied per-cell -- a single comma-
# separated multi-area Range only actually formats its first area.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.674.60"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.36%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.106.19"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.92%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "527.53"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.80"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.91%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.96%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.35"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.48%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.60%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.637.32"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.94"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +5.06%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "58.668.26"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.15"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.089.66"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.95"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.16%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "342.79"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.53%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.507"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.96"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.26%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.52%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0₃0920"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.67"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +3.98%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.26"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.16%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.89%  "
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.21"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +3.49%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.07"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "154.35"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.44%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.09"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "27.14"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.49%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +5.28%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.45%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.142.18"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.67%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.92"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.78%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.674"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.34%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +7.08%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.286.04"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.55%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "20.98"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +4.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.970"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.80%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.03%  "
$ws.Range("B50").Value = "SuiNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.751"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +9.27%  "
$ws.Range("B51").Value = "Bittensor"
$ws.Range("C51").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "268.23"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +7.64%  "
